$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '64.587.22'
$ws.Range("E2").Value = '  -0.98%  '

Set-TextValue "D3" '3.440.43'
$ws.Range("E3").Value = '  -1.32%  '

$ws.Range("E4").Value = '  +0.06%  '

Set-TextValue "D5" '572.32'
$ws.Range("E5").Value = '  -1.36%  '

Set-TextValue "D6" '158.46'
$ws.Range("E6").Value = '  -2.73%  '

$ws.Range("E7").Value = '  +0.04%  '

Set-TextValue "D8" '3.438.91'
$ws.Range("E8").Value = '  -1.31%  '

Set-TextValue "D9" '0.571'
$ws.Range("E9").Value = '  -7.21%  '

Set-TextValue "D10" '7.20'
$ws.Range("E10").Value = '  -0.73%  '

$ws.Range("E11").Value = '  -3.81%  '

Set-TextValue "D12" '0.436'
$ws.Range("E12").Value = '  -2.63%  '

Set-TextValue "D13" '4.034.93'
$ws.Range("E13").Value = '  -1.24%  '

$ws.Range("E14").Value = '  -0.55%  '

Set-TextValue "D15" '27.43'
$ws.Range("E15").Value = '  -4.52%  '

Set-TextValue "D16" '0.0000173'
$ws.Range("E16").Value = '  -10.64%  '

Set-TextValue "D17" '64.683.86'
$ws.Range("E17").Value = '  -0.82%  '

Set-TextValue "D18" '3.439.70'
$ws.Range("E18").Value = '  -1.58%  '

Set-TextValue "D19" '6.14'
$ws.Range("E19").Value = '  -5.17%  '

Set-TextValue "D20" '13.66'
$ws.Range("E20").Value = '  -5.06%  '

Set-TextValue "D21" '377.10'
$ws.Range("E21").Value = '  -1.67%  '

Set-TextValue "D22" '7.85'
$ws.Range("E22").Value = '  -4.54%  '

$ws.Range("E23").Value = '  -0.13%  '

Set-TextValue "D24" '72.02'
$ws.Range("E24").Value = '  -0.84%  '

Set-TextValue "D25" '0.531'
$ws.Range("E25").Value = '  -4.29%  '

Set-TextValue "D26" '0.0000119'
$ws.Range("E26").Value = '  -1.32%  '

Set-TextValue "D27" '9.88'
$ws.Range("E27").Value = '  -1.59%  '

$ws.Range("E28").Value = '  -0.15%  '

Set-TextValue "D29" '1.00'
$ws.Range("E29").Value = '  +0.07%  '

Set-TextValue "D30" '1.43'
$ws.Range("E30").Value = '  -7.25%  '

Set-TextValue "D31" '6.03'
$ws.Range("E31").Value = '  -2.45%  '

Set-TextValue "D32" '2.00'
$ws.Range("E32").Value = '  -2.62%  '

Set-TextValue "D33" '23.12'
$ws.Range("E33").Value = '  -2.44%  '

Set-TextValue "D34" '6.95'
$ws.Range("E34").Value = '  -3.44%  '

$ws.Range("E35").Value = '  -5.10%  '

Set-TextValue "D36" '160.83'
$ws.Range("E36").Value = '  -0.83%  '

Set-TextValue "D37" '1.87'
$ws.Range("E37").Value = '  -3.10%  '

Set-TextValue "D38" '2.883.48'
$ws.Range("E38").Value = '  -4.08%  '

Set-TextValue "D39" '0.0743'
$ws.Range("E39").Value = '  -4.87%  '

Set-TextValue "D40" '25.98'
$ws.Range("E40").Value = '  -3.49%  '

Set-TextValue "D41" '43.00'
$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("E42").Value = '  +0.77%  '

Set-TextValue "D43" '4.49'
$ws.Range("E43").Value = '  -1.95%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D44" '6.43'
$ws.Range("E44").Value = '  -5.57%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D45" '25.89'
$ws.Range("E45").Value = '  -0.35%  '

$ws.Range("E46").Value = '  -4.23%  '

Set-TextValue "D47" '2.36'
$ws.Range("E47").Value = '  +7.92%  '

Set-TextValue "D48" '320.76'
$ws.Range("E48").Value = '  +0.71%  '

Set-TextValue "D49" '1.07'
$ws.Range("E49").Value = '  -3.82%  '

Set-TextValue "D50" '6.44'
$ws.Range("E50").Value = '  -4.42%  '

Set-TextValue "D51" '0.839'
$ws.Range("E51").Value = '  -4.28%  '
